$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$names = @(
    "Danielle Seunarine",
    "Jochelle Greaves",
    "Salma Kaloo",
    "Celina Mahabir",
    "Lilian Soogrim",
    "Asmita Nankissoon",
    "Kady Seecharan",
    "Shivanna Sookdeo",
    "Janelle Raghoo",
    "Analeyah Ryan"
)

$row = 5
foreach ($name in $names) {
    $ws.Range("A$row").Value = $name
    $row++
}

$ws.Columns.Item(1).ColumnWidth = 16.5546875

$ws.Range("A15").Select()
